# Generate Report for Handoff
#
# The two tracked files (ae028673-...md and cd54ed06-...md) swap row
# positions across all three sheets, and the ae028673 file's status moves
# from "Handed back: in sync with en-US" to "Ready for handoff" (with a new
# HO-xliff handoff timestamp and a "stale handback" error message recorded
# on the language sheets). The "Error Detail" column is widened to fit the
# longer message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": rows swap (cd54ed06 now row 2, ae028673 now row 3) and
# ae028673's status/date (now in row 3) becomes "Ready for handoff".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.md"
$wsOverview.Range("A3").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.md"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-28 06:53:06"

# Hyperlinks B2/B3 keep pointing at their original targets (rId2 -> ae028673
# URL, rId3 -> cd54ed06 URL) but now display the swapped file name, matching
# the new row contents. Rebuild the hyperlink list to flip the display text
# while leaving each relationship's target untouched.
$wsOverview.Range("A1").Hyperlinks.Delete() | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d71a41437d89b354f62bbc01cb4a28f6932c48dc/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md", "", "", "e2e\cd54ed06-4bee-4486-a1f0-1dc02011ca95.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d71a41437d89b354f62bbc01cb4a28f6932c48dc/e2e/cd54ed06-4bee-4486-a1f0-1dc02011ca95.md", "", "", "e2e\ae028673-f9a0-4771-8ab2-37d5a1c491fe.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": same row swap; ae028673 (now row 3) is "Ready for handoff"
# with a new Latest Handoff File/Datetime and an Error Detail message.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("G2").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.12ffa7c52420325959c5e575bf27d9b8c17ed3d5.zh-cn.xlf"
$wsZh.Range("J2").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.12ffa7c52420325959c5e575bf27d9b8c17ed3d5.zh-cn.xlf"
$wsZh.Range("G3").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.e8be8300cd87d911dc40d54f726a2dda27c158dc.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-28 06:52:58"
$wsZh.Range("J3").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.e8be8300cd87d911dc40d54f726a2dda27c158dc.zh-cn.xlf"
$wsZh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d71a41437d89b354f62bbc01cb4a28f6932c48dc/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc0e7f58f95a424f1c872f23ff4cabe3e4142040/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md."

# Rebuild the A2/I2/A3/I3 hyperlinks the same way as on Overview: each rId
# keeps its original target, only the displayed file name flips.
$wsZh.Range("A1").Hyperlinks.Delete() | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d71a41437d89b354f62bbc01cb4a28f6932c48dc/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md", "", "", "cd54ed06-4bee-4486-a1f0-1dc02011ca95.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f75a0f86ed382d50c3fe36f2b311bdebe3545c8d/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md", "", "", "cd54ed06-4bee-4486-a1f0-1dc02011ca95.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d71a41437d89b354f62bbc01cb4a28f6932c48dc/e2e/cd54ed06-4bee-4486-a1f0-1dc02011ca95.md", "", "", "ae028673-f9a0-4771-8ab2-37d5a1c491fe.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f75a0f86ed382d50c3fe36f2b311bdebe3545c8d/e2e/cd54ed06-4bee-4486-a1f0-1dc02011ca95.md", "", "", "ae028673-f9a0-4771-8ab2-37d5a1c491fe.md") | Out-Null

# Widen the "Error Detail" column (P, 16th) to fit the new message.
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# Sheet "de-de": identical pattern to "zh-cn".
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("G2").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.12ffa7c52420325959c5e575bf27d9b8c17ed3d5.de-de.xlf"
$wsDe.Range("J2").Value = "cd54ed06-4bee-4486-a1f0-1dc02011ca95.12ffa7c52420325959c5e575bf27d9b8c17ed3d5.de-de.xlf"
$wsDe.Range("G3").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.e8be8300cd87d911dc40d54f726a2dda27c158dc.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-28 06:53:06"
$wsDe.Range("J3").Value = "ae028673-f9a0-4771-8ab2-37d5a1c491fe.e8be8300cd87d911dc40d54f726a2dda27c158dc.de-de.xlf"
$wsDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d71a41437d89b354f62bbc01cb4a28f6932c48dc/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc0e7f58f95a424f1c872f23ff4cabe3e4142040/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md."

$wsDe.Range("A1").Hyperlinks.Delete() | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d71a41437d89b354f62bbc01cb4a28f6932c48dc/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md", "", "", "cd54ed06-4bee-4486-a1f0-1dc02011ca95.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a08cff03d9a5981a037fee6bf96e0ee4692dad22/e2e/ae028673-f9a0-4771-8ab2-37d5a1c491fe.md", "", "", "cd54ed06-4bee-4486-a1f0-1dc02011ca95.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d71a41437d89b354f62bbc01cb4a28f6932c48dc/e2e/cd54ed06-4bee-4486-a1f0-1dc02011ca95.md", "", "", "ae028673-f9a0-4771-8ab2-37d5a1c491fe.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a08cff03d9a5981a037fee6bf96e0ee4692dad22/e2e/cd54ed06-4bee-4486-a1f0-1dc02011ca95.md", "", "", "ae028673-f9a0-4771-8ab2-37d5a1c491fe.md") | Out-Null

$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
